# Auto-generated edit script to update cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells remain stored as text, matching the source
# data where prices use '.' as a thousands-like separator rather than a
# decimal point (e.g. "28.586.94"), which Excel would otherwise reinterpret
# as a number and mangle.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '28.586.94'
$ws.Cells.Item(2, 5).Value = '  +2.34%  '
$ws.Cells.Item(3, 4).Value = '1.913.88'
$ws.Cells.Item(4, 5).Value = '  -0.20%  '
$ws.Cells.Item(5, 4).Value = '315.69'
$ws.Cells.Item(5, 5).Value = '  +1.70%  '
$ws.Cells.Item(6, 4).Value = '1.000'
$ws.Cells.Item(6, 5).Value = '  -0.11%  '
$ws.Cells.Item(7, 4).Value = '0.5145'
$ws.Cells.Item(7, 5).Value = '  +2.99%  '
$ws.Cells.Item(8, 4).Value = '0.3979'
$ws.Cells.Item(8, 5).Value = '  +1.81%  '
$ws.Cells.Item(9, 4).Value = '0.09813'
$ws.Cells.Item(9, 5).Value = '  +0.97%  '
$ws.Cells.Item(10, 4).Value = '1.162'
$ws.Cells.Item(10, 5).Value = '  +5.70%  '
$ws.Cells.Item(11, 4).Value = '42.40'
$ws.Cells.Item(11, 5).Value = '  +3.61%  '
$ws.Cells.Item(12, 4).Value = '6.555'
$ws.Cells.Item(12, 5).Value = '  +2.46%  '
$ws.Cells.Item(13, 4).Value = '21.31'
$ws.Cells.Item(13, 5).Value = '  +4.26%  '
$ws.Cells.Item(14, 4).Value = '1.915.45'
$ws.Cells.Item(14, 5).Value = '  +5.49%  '
$ws.Cells.Item(15, 4).Value = '7.596'
$ws.Cells.Item(15, 5).Value = '  +4.67%  '
$ws.Cells.Item(16, 4).Value = '1.000'
$ws.Cells.Item(16, 5).Value = '  -0.19%  '
$ws.Cells.Item(17, 4).Value = '0.00001145'
$ws.Cells.Item(17, 5).Value = '  +0.75%  '
$ws.Cells.Item(18, 4).Value = '94.11'
$ws.Cells.Item(18, 5).Value = '  +1.92%  '
$ws.Cells.Item(19, 4).Value = '0.06661'
$ws.Cells.Item(19, 5).Value = '  +0.27%  '
$ws.Cells.Item(20, 4).Value = '18.22'
$ws.Cells.Item(20, 5).Value = '  +6.27%  '
$ws.Cells.Item(21, 5).Value = '  -0.13%  '
$ws.Cells.Item(22, 4).Value = '6.332'
$ws.Cells.Item(22, 5).Value = '  +7.00%  '
$ws.Cells.Item(23, 4).Value = '28.643.87'
$ws.Cells.Item(23, 5).Value = '  +2.31%  '
$ws.Cells.Item(24, 5).Value = '  +3.71%  '
$ws.Cells.Item(25, 4).Value = '2.296'
$ws.Cells.Item(25, 5).Value = '  +1.89%  '
$ws.Cells.Item(26, 4).Value = '2.716'
$ws.Cells.Item(26, 5).Value = '  +13.87%  '
$ws.Cells.Item(27, 4).Value = '2.133.29'
$ws.Cells.Item(27, 5).Value = '  +5.57%  '
$ws.Cells.Item(28, 4).Value = '21.32'
$ws.Cells.Item(28, 5).Value = '  +3.80%  '
$ws.Cells.Item(29, 4).Value = '159.84'
$ws.Cells.Item(29, 5).Value = '  +0.74%  '
$ws.Cells.Item(30, 4).Value = '129.09'
$ws.Cells.Item(30, 5).Value = '  +1.91%  '
$ws.Cells.Item(31, 4).Value = '1.106'
$ws.Cells.Item(31, 5).Value = '  +7.29%  '
$ws.Cells.Item(32, 4).Value = '0.1084'
$ws.Cells.Item(32, 5).Value = '  +1.83%  '
$ws.Cells.Item(33, 4).Value = '5.786'
$ws.Cells.Item(33, 5).Value = '  +4.28%  '
$ws.Cells.Item(34, 4).Value = '3.639'
$ws.Cells.Item(34, 5).Value = '  +1.48%  '
$ws.Cells.Item(35, 4).Value = '9.900'
$ws.Cells.Item(35, 5).Value = '  +11.22%  '
$ws.Cells.Item(36, 4).Value = '0.06821'
$ws.Cells.Item(36, 5).Value = '  +1.40%  '
$ws.Cells.Item(37, 5).Value = '  +5.19%  '
$ws.Cells.Item(38, 4).Value = '1.269'
$ws.Cells.Item(38, 5).Value = '  +7.77%  '
$ws.Cells.Item(39, 4).Value = '0.2238'
$ws.Cells.Item(39, 5).Value = '  +4.63%  '
$ws.Cells.Item(40, 5).Value = '  +6.29%  '
$ws.Cells.Item(41, 4).Value = '5.125'
$ws.Cells.Item(41, 5).Value = '  +3.90%  '
$ws.Cells.Item(42, 4).Value = '0.6463'
$ws.Cells.Item(42, 5).Value = '  +4.82%  '
$ws.Cells.Item(43, 4).Value = '1.196'
$ws.Cells.Item(43, 5).Value = '  +2.40%  '
$ws.Cells.Item(44, 5).Value = '  -0.06%  '
$ws.Cells.Item(45, 5).Value = '  +4.07%  '
$ws.Cells.Item(46, 4).Value = '0.6110'
$ws.Cells.Item(46, 5).Value = '  +3.83%  '
$ws.Cells.Item(47, 4).Value = '3.807'
$ws.Cells.Item(47, 5).Value = '  +3.06%  '
$ws.Cells.Item(48, 4).Value = '1.280'
$ws.Cells.Item(48, 5).Value = '  -0.32%  '
$ws.Cells.Item(49, 5).Value = '  +5.83%  '
$ws.Cells.Item(50, 4).Value = '125.78'
$ws.Cells.Item(50, 5).Value = '  +1.85%  '
$ws.Cells.Item(51, 4).Value = '1.217'
$ws.Cells.Item(51, 5).Value = '  +3.45%  '
